# Apply the "Test Cases" sheet updates described in the commit:
#  - Correct the TestCase_B17 description (was copy/pasted as the "ALL" content
#    type / documents wording; should reference "articles" / "ARTICLES content type")
#  - Mark TestCase_B41/B42/B43 results as SKIP (only the last case, B44, stays PASS)
#  - Update the active selection left after this edit to D4 (clearing the old
#    left-scrolled viewport)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C18").Value = "Verify that user is able to sort the articles by TIMES CITED field in ARTICLES content type"

$ws.Range("E42").Value = "SKIP"
$ws.Range("E43").Value = "SKIP"
$ws.Range("E44").Value = "SKIP"

# Leave the selection/viewport where the author left it after editing.
$ws.Range("D4").Select()
